# Update the "想去人数" (interested-count) figures in the F column of the
# "展览" (Exhibition) and "全部类型" (All types) sheets to reflect a newer
# scrape of the underlying data.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 400
$ws1.Range("F4").Value = 5050
$ws1.Range("F5").Value = 40
$ws1.Range("F6").Value = 39
$ws1.Range("F8").Value = 498

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 400
$ws4.Range("F4").Value = 5050
$ws4.Range("F6").Value = 40
$ws4.Range("F7").Value = 39
$ws4.Range("F10").Value = 498
